# Scripts/Industry/Results_per_Country/2050_MT.xlsx — "Updated Results with corrected code"
#
# The sheet lists demand-by-energy-carrier rows. The corrected run:
#   1. stops reporting a (spurious) 0 for Hydrogen's "Non-metallic minerals" cell (D3)
#   2. realizes row 7 was mislabeled "Other" when it is actually "Biogas"
#   3. appends the real trailing "Other" row (row 8) that row 7 used to represent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) D3 is no longer a computed 0 under the corrected code - blank it out.
$ws.Range("D3").ClearContents()

# 2) Row 7's category label was "Other" but should read "Biogas".
$ws.Range("A7").Value = "Biogas"

# 3) Add the new trailing "Other" row 8, carrying over row 7's old
#    row-label formatting (bold, bordered, centered/top-aligned) and
#    its blank B/C cells + numeric 0 in D.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 0
